$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in the test case numbers for Day 8 (Written / Execution / Review)
$ws.Range("C49").Value = 696
$ws.Range("C50").Value = 882
$ws.Range("C51").Value = 615

# Update the view's scroll position and active selection
$excel.ActiveWindow.ScrollRow = 43
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("C51").Select()
